$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# New identifiers replacing the old ones
$oldId1 = "87125e70-4812-4d7a-bab9-591f8a17caf5"
$newId1 = "33d99b1b-f6df-4c6b-946b-effcbb21a229"
$oldId2 = "f9fde33b-ca70-47b4-998d-c05cc45437ce"
$newId2 = "ffff809b08e7-ae90-45ee-93b3-582579312b3e"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-13 17:28:55"
$newHandoffDate = "2016-08-13 17:28:47"
$newHandbackDate = "0001-01-01 00:00:00"

$newHandoffFile = "$newId1.660067b065ccdcb6ad2a13e5b5e2f7606ee8ed26.zh-cn.xlf"
$newHandoffFileDe = "$newId1.660067b065ccdcb6ad2a13e5b5e2f7606ee8ed26.de-de.xlf"

# ---- Overview sheet ----
$overview.Range("A2").Value = "$newId1.md"
$overview.Range("B2").Value = "e2e\$newId1.md"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = $newHoDate

$overview.Range("A3").Value = "$newId2.md"
$overview.Range("B3").Value = "e2e\$newId2.md"
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newHoDate

$overview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newId1.md"
$overview.Hyperlinks.Item(2).TextToDisplay = "e2e\$newId2.md"

# ---- zh-cn sheet ----
$zhcn.Range("A2").Value = "$newId1.md"
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("G2").Value = $newHandoffFile
$zhcn.Range("H2").Value = $newHandoffDate
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = $newHandbackDate

$zhcn.Range("A3").Value = "$newId2.md"
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("F3").Value = "True"
$zhcn.Range("G3").Value = $newHandoffFile
$zhcn.Range("H3").Value = $newHandoffDate
$zhcn.Range("I3").Value = ""
$zhcn.Range("J3").Value = ""
$zhcn.Range("K3").Value = $newHandbackDate

$zhcn.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$zhcn.Hyperlinks.Item(2).TextToDisplay = "$newId2.md"

# ---- de-de sheet ----
$dede.Range("A2").Value = "$newId1.md"
$dede.Range("C2").Value = $newStatus
$dede.Range("G2").Value = $newHandoffFileDe
$dede.Range("H2").Value = $newHoDate
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = $newHandbackDate

$dede.Range("A3").Value = "$newId2.md"
$dede.Range("C3").Value = $newStatus
$dede.Range("F3").Value = "True"
$dede.Range("G3").Value = $newHandoffFileDe
$dede.Range("H3").Value = $newHoDate
$dede.Range("I3").Value = ""
$dede.Range("J3").Value = ""
$dede.Range("K3").Value = $newHandbackDate

$dede.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$dede.Hyperlinks.Item(2).TextToDisplay = "$newId2.md"
